$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "54.584.77"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.276.72"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "2.294.72"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "2.684.22"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").Value = "54.646.59"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000130"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "2.292.57"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "306.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.909"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.54%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "249.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.23%  "
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
